# Add a "Cat" / "Something" pair of columns (D/E) to the existing formula
# test sheet, mirroring the style of the existing Sin/Cos helper columns:
#   - D1/E1 get new header labels
#   - D2 seeds a value, D3:D16 build a shared CONCAT formula chain
#   - column D gets a wider custom width so the longer strings are visible
#   - leave the selection where the user's last click landed (E17)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Cat"
$ws.Range("E1").Value = "Something"

$ws.Range("D2").Value = 1
$ws.Range("D3").Formula = "=CONCAT(D2,""a"")"
$ws.Range("D4:D16").Formula = "=CONCAT(D3,""a"")"

$ws.Range("D1").ColumnWidth = 17.3

$ws.Range("E17").Select() | Out-Null
